# Updates the cryptos price/volume table to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.278.56'
$ws.Range('E2').Value = '  -0.16%  '

# Row 3
$ws.Range('D3').Value = '1.929.98'
$ws.Range('E3').Value = '  -0.24%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.91'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.22%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7166'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.51%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3189'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.73%  '

# Row 9
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.73'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.37%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07107'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.30%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7912'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.03%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07983'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.53%  '

# Row 13
$ws.Range('D13').Value = '1.926.91'
$ws.Range('E13').Value = '  -0.39%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.389'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.71%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.91'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.12%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.69'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.96%  '

# Row 17
$ws.Range('D17').Value = '30.266.42'
$ws.Range('E17').Value = '  -0.20%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '255.64'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008063'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.06%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.770'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.48%  '

# Row 21
$ws.Range('D21').Value = '2.179.13'
$ws.Range('E21').Value = '  -0.39%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9995'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.04%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.11%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.829'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.65%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.540'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.78%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.20'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.13%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.12'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.82%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.263'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.84%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1268'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.73%  '

# Row 30
$ws.Range('E30').Value = '  +0.79%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.526'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.18%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.399'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.97%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.132'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.37%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05139'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.04%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.274'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7456'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.73%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.761'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.01%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01962'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.77%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.796'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.33%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.71'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.65%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.365'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.45%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4513'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.20%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.988'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.58%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8469'
$ws.Range('D44').ClearFormats()

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9993'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('E46').Value = '  -2.04%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.780'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.19%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.419'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.64%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.66'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.41%  '

# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '949.14'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +9.91%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06104'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.14%  '

